# Weekly data refresh for the Mango price series:
# insert a new observation row at row 122 (pushing the existing
# rows 122-158 down to 123-159) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 122..158 down to 123..159, leaving a blank row 122.
$ws.Rows.Item(122).Insert()

# Fill in the new row 122 with the latest weekly price observation.
$ws.Cells.Item(122, 1).Value  = 7
$ws.Cells.Item(122, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value  = "Ñuble"
$ws.Cells.Item(122, 4).Value  = 45119
$ws.Cells.Item(122, 5).Value  = 16
$ws.Cells.Item(122, 6).Value  = "Fruta"
$ws.Cells.Item(122, 7).Value  = 100108
$ws.Cells.Item(122, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(122, 9).Value  = 100108002
$ws.Cells.Item(122, 10).Value = "Mango"
$ws.Cells.Item(122, 11).Value = "Sin especificar"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 30
$ws.Cells.Item(122, 14).Value = 9000
$ws.Cells.Item(122, 15).Value = 9000
$ws.Cells.Item(122, 16).Value = 9000
$ws.Cells.Item(122, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(122, 18).Value = "Brasil"
$ws.Cells.Item(122, 19).Value = 2250
$ws.Cells.Item(122, 20).Value = 4
